# Updated cryptos list with latest prices and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.248.17'
$ws.Range("E2").Value = '  -4.03%  '

# Row 3
$ws.Range("D3").Value = '2.954.08'
$ws.Range("E3").Value = '  -6.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.86'
$ws.Range("E5").Value = '  -5.07%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.24'
$ws.Range("E6").Value = '  -6.74%  '

# Row 7
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  -3.82%  '

# Row 9
$ws.Range("D9").Value = '2.958.60'
$ws.Range("E9").Value = '  -5.91%  '

# Row 10
$ws.Range("E10").Value = '  -4.24%  '

# Row 11
$ws.Range("E11").Value = '  -8.19%  '

# Row 12
$ws.Range("E12").Value = '  -4.99%  '

# Row 13
$ws.Range("D13").Value = '3.470.82'
$ws.Range("E13").Value = '  -5.93%  '

# Row 14
$ws.Range("E14").Value = '  -2.43%  '

# Row 15
$ws.Range("D15").Value = '61.247.91'
$ws.Range("E15").Value = '  -4.40%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.52'
$ws.Range("E16").Value = '  -6.53%  '

# Row 17
$ws.Range("D17").Value = '2.961.12'
$ws.Range("E17").Value = '  -6.03%  '

# Row 18
$ws.Range("E18").Value = '  -6.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.08'
$ws.Range("E19").Value = '  -2.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '380.17'
$ws.Range("E20").Value = '  -5.76%  '

# Row 21
$ws.Range("E21").Value = '  -5.83%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").Value = '  -5.91%  '

# Row 23
$ws.Range("E23").Value = '  +0.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.76'
$ws.Range("E24").Value = '  -5.18%  '

# Row 25
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.078.35'
$ws.Range("E25").Value = '  -7.20%  '

# Row 26
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.465'
$ws.Range("E26").Value = '  -3.48%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.184'
$ws.Range("E27").Value = '  -7.24%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0917'
$ws.Range("E29").Value = '  -9.92%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.23'
$ws.Range("E30").Value = '  -6.12%  '

# Row 31
$ws.Range("E31").Value = '  -0.03%  '

# Row 32
$ws.Range("E32").Value = '  -5.69%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.18'
$ws.Range("E33").Value = '  -4.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '157.00'
$ws.Range("E34").Value = '  +0.77%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.93'
$ws.Range("E35").Value = '  -5.75%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.53'
$ws.Range("E36").Value = '  -6.47%  '

# Row 37
$ws.Range("E37").Value = '  -6.09%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.26'
$ws.Range("E38").Value = '  -5.31%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '2.417.34'
$ws.Range("E39").Value = '  -9.96%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.52'
$ws.Range("E40").Value = '  -9.40%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.88'
$ws.Range("E41").Value = '  -4.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.03'
$ws.Range("E42").Value = '  -7.15%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '36.70'
$ws.Range("E43").Value = '  -4.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.659'
$ws.Range("E44").Value = '  -4.94%  '

# Row 45
$ws.Range("E45").Value = '  -4.72%  '

# Row 46
$ws.Range("E46").Value = '  -0.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0243'
$ws.Range("E47").Value = '  -5.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.84'
$ws.Range("E48").Value = '  -9.51%  '

# Row 49
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.49'
$ws.Range("E49").Value = '  +0.24%  '

# Row 50
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0942'
$ws.Range("E50").Value = '  -3.61%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.44'
$ws.Range("E51").Value = '  -8.31%  '
